$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
try {
    $ws.Cells.Item(8,6).Value = "Organisation"
    Write-Output "cells item ok"
} catch {
    Write-Output ("cells item err: " + $_.Exception.Message)
}
try {
    $ws.Columns.Item(6).Cells.Item(1,1).Value = "x"
    Write-Output "columns ok"
} catch {
    Write-Output ("columns err: " + $_.Exception.Message)
}
